$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily auto-push: insert the newest timestamp record as a new row at 835,
# which pushes all following rows (formerly 835..876) down to 836..877.
$ws.Rows.Item(835).Insert()

# Column A holds dates stored as plain text (e.g. "2026/12/29"), so force a
# text format while assigning the value to stop Excel from auto-converting
# the date-shaped string into a date serial, then clear the format override
# again so the new cell ends up with the same default (unstyled) formatting
# as every other data row.
$ws.Range("A835").NumberFormat = "@"
$ws.Range("A835").Value = "2026/02/20"
$ws.Range("A835").ClearFormats()

$ws.Range("B835").Value = "金"
$ws.Range("C835").Value = 10
$ws.Range("D835").Value = 201
